{"js": "// UART protocol doc: \"V\u1edbi 1 byte ta c\u00f3 th\u1ec3 m\u00e3 h\u00f3a \u0111\u01b0\u1ee3c 128 lo\u1ea1i c\u00e2u l\u1ec7nh\n// kh\u00e1c nhau.\" -> 1 byte actually encodes 256 distinct values, so the\n// figure is corrected from 128 to 256.\nconst body = context.document.body;\n\nlet results = body.search(\"128 lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau\", { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  // Fallback: narrower search in case surrounding text ever shifts.\n  results = body.search(\"128\", { matchCase: true, matchWholeWord: true });\n  results.load(\"items/text\");\n  await context.sync();\n}\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '128' (lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau) to update.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  const item = results.items[i];\n  const newText = item.text.replace(\"128\", \"256\");\n  item.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# UART protocol doc: \"V\u1edbi 1 byte ta c\u00f3 th\u1ec3 m\u00e3 h\u00f3a \u0111\u01b0\u1ee3c 128 lo\u1ea1i c\u00e2u l\u1ec7nh\n# kh\u00e1c nhau.\" -> 1 byte actually encodes 256 distinct values, so the\n# figure is corrected from 128 to 256.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"128 lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"256 lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$found = $find.Execute(\"128 lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau\", $false, $false, $false, $false, $false, $true, 1, $false, \"256 lo\u1ea1i c\u00e2u l\u1ec7nh kh\u00e1c nhau\", 2)\n\nif (-not $found) {\n    # Fallback: narrower search in case surrounding text ever shifts.\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Text = \"128\"\n    $find2.Replacement.ClearFormatting()\n    $find2.Replacement.Text = \"256\"\n    $find2.Execute(\"128\", $false, $true, $false, $false, $false, $true, 1, $false, \"256\", 2) | Out-Null\n}\n"}
